$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3622736666666667
$ws.Range("N2").Value = 1.086821
$ws.Range("O2").Value = 0.02355627613074877
$ws.Range("P2").Value = 0.02355627613074877
$ws.Range("Q2").Value = 0.05792768005788889
$ws.Range("R2").Value = 0.521349120521
$ws.Range("S2").Value = 0.0006237164703157517
$ws.Range("T2").Value = 0.0006237164703157518

# Row 3
$ws.Range("O3").Value = 0.8197376278620713
$ws.Range("P3").Value = 0.8197376278620713
$ws.Range("S3").Value = 0.02170478292057984
$ws.Range("T3").Value = 0.02170478292057984

# Row 4
$ws.Range("M4").Value = 2.409994333333333
$ws.Range("N4").Value = 7.229983
$ws.Range("O4").Value = 0.15670609600718
$ws.Range("P4").Value = 0.15670609600718
$ws.Range("Q4").Value = 0.3853588972314444
$ws.Range("R4").Value = 3.468230075083
$ws.Range("S4").Value = 0.004149220043781716
$ws.Range("T4").Value = 0.004149220043781718

# Row 5
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.3622736666666667
$ws.Range("N5").Value = 1.086821
$ws.Range("O5").Value = 0.02355627613074877
$ws.Range("P5").Value = 0.02355627613074877
$ws.Range("Q5").Value = 2.129861951930667
$ws.Range("R5").Value = 19.168757567376
$ws.Range("S5").Value = 0.02293255966043302
$ws.Range("T5").Value = 0.02293255966043302

# Row 6
$ws.Range("O6").Value = 0.8197376278620713
$ws.Range("P6").Value = 0.8197376278620713
$ws.Range("Q6").Value = 74.11731695020799
$ws.Range("R6").Value = 667.0558525518719
$ws.Range("S6").Value = 0.7980328449414914
$ws.Range("T6").Value = 0.7980328449414915

# Row 7
$ws.Range("M7").Value = 2.409994333333333
$ws.Range("N7").Value = 7.229983
$ws.Range("O7").Value = 0.15670609600718
$ws.Range("P7").Value = 0.15670609600718
$ws.Range("Q7").Value = 14.16872300480533
$ws.Range("R7").Value = 127.518507043248
$ws.Range("S7").Value = 0.1525568759633982
$ws.Range("T7").Value = 0.1525568759633983
